$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.23"
$ws.Range("E2").Value = "'2.85%"
$ws.Range("D3").Value = "'41.24"
$ws.Range("E3").Value = "'2.45%"
$ws.Range("D4").Value = "'5.006"
$ws.Range("E4").Value = "'-0.71%"
$ws.Range("D5").Value = "'0.07518"
$ws.Range("E5").Value = "'2.99%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.362"
$ws.Range("E6").Value = "'1.79%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.568"
$ws.Range("E7").Value = "'2.84%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9295"
$ws.Range("E8").Value = "'1.08%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.408"
$ws.Range("E9").Value = "'0.46%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1218"
$ws.Range("E10").Value = "'2.87%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1840"
$ws.Range("E11").Value = "'6.29%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08865"
$ws.Range("E12").Value = "'1.92%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04101"
$ws.Range("E13").Value = "'-1.73%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1053"
$ws.Range("E14").Value = "'0.06%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001285"
$ws.Range("E15").Value = "'1.39%"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04045"
$ws.Range("E16").Value = "'4.63%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005955"
$ws.Range("E17").Value = "'-0.10%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.344"
$ws.Range("E18").Value = "'-1.61%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3288"
$ws.Range("E19").Value = "'-0.07%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'7.977"
$ws.Range("E20").Value = "'1.57%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1418"
$ws.Range("E21").Value = "'5.76%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2965"
$ws.Range("E22").Value = "'2.85%"
$ws.Range("E23").Value = "'-0.26%"
$ws.Range("D24").Value = "'0.003885"
$ws.Range("E24").Value = "'1.24%"
$ws.Range("E25").Value = "'-4.05%"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("D38").Value = "'0.02421"
$ws.Range("E38").Value = "'4.11%"
$ws.Range("D39").Value = "'0.05212"
$ws.Range("E39").Value = "'4.53%"
$ws.Range("D40").Value = "'0.005903"
$ws.Range("E40").Value = "'-8.65%"
$ws.Range("D41").Value = "'0.007793"
$ws.Range("E41").Value = "'1.68%"
$ws.Range("D42").Value = "'0.1329"
$ws.Range("E42").Value = "'4.10%"
$ws.Range("D43").Value = "'0.007366"
$ws.Range("E43").Value = "'0.17%"
$ws.Range("D44").Value = "'0.007829"
$ws.Range("E44").Value = "'10.83%"
$ws.Range("D45").Value = "'0.2968"
$ws.Range("E45").Value = "'-5.24%"
$ws.Range("D46").Value = "'0.00006324"
$ws.Range("E46").Value = "'-1.85%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("D48").Value = "'0.04517"
$ws.Range("E48").Value = "'429.34%"
$ws.Range("D49").Value = "'0.004196"
$ws.Range("E49").Value = "'-0.12%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("D51").Value = "'0.0001998"
